$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-12-29 Friday" "2023-12-30 Saturday"

Replace-Text "88×56=" "63×54="
Replace-Text "86×91=" "80×16="
Replace-Text "97×46=" "34×63="
Replace-Text "25×33=" "77×19="
Replace-Text "23×78=" "26×86="

Replace-Text "30×37=" "66×17="
Replace-Text "66×23=" "70×88="
Replace-Text "25×69=" "32×16="
Replace-Text "21×12=" "27×92="
Replace-Text "33×75=" "46×34="

Replace-Text "85×59=" "69×96="
Replace-Text "79×39=" "90×62="
Replace-Text "22×84=" "74×90="
Replace-Text "94×29=" "81×42="
Replace-Text "57×62=" "63×89="

Replace-Text "90×80=" "96×63="
Replace-Text "94×97=" "28×16="
Replace-Text "14×57=" "99×54="
Replace-Text "28×23=" "64×61="
Replace-Text "18×46=" "69×73="

Replace-Text "29×50=" "78×99="
Replace-Text "45×26=" "62×38="
Replace-Text "66×75=" "37×93="
Replace-Text "47×98=" "53×85="
Replace-Text "57×26=" "93×56="
